$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "user_id" column before the old column B (id_catarticle),
# shifting the rest of the table (old B..M) one column to the right (new C..N).
$ws.Columns("B").Insert()

# Stash the two distinct formats used by the hyperlink cells (plain-looking
# "Hyperlink" style without border, and with border) before anything else
# touches them, so we can restore them after Hyperlinks.Add() below forces
# its own (underlined / colored) look on the cell it is attached to.
$ws.Range("H2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("Z2").PasteSpecial(-4122)

# New header + data for the inserted column.
$ws.Range("B1").Value = "user_id"
$ws.Range("B2:B6").Value = 1

# The column insert shifted the actual hyperlinked cell content from the old
# column G to the new column H, but left the <hyperlinks> bookkeeping
# pointing at the stale G column - recreate each hyperlink against its new
# column H location (same addresses/order as before the insert).
$ws.Range("G2:G6").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H6"), "https://aminama.com/bahan-kain-jersey/")
$ws.Hyperlinks.Add($ws.Range("H5"), "https://aminama.com/bahan-kain-jersey/")
$ws.Hyperlinks.Add($ws.Range("H2"), "https://tshirtbar.id/perbedaan-cotton-combed/")
$ws.Hyperlinks.Add($ws.Range("H4"), "https://aminama.com/bahan-fleece/")
$ws.Hyperlinks.Add($ws.Range("H3"), "https://aminama.com/bahan-fleece/")

# Restore original formatting now that hyperlinks point at the right cells.
$ws.Range("Z1").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)

$ws.Range("Z2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H5").PasteSpecial(-4122)

$ws.Range("Z1:Z2").Clear()

# Update selection to match the new layout.
$ws.Range("B7").Select()
